$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D2 (Basic progress) from 6 to 11 - best solution added
$ws.Range("D2").Value = 11

# Set D3:D16 (rest of the "Basic" rows) to 0 so the SUM(D2:D16) range is fully populated
for ($r = 3; $r -le 16; $r++) {
    $ws.Cells.Item($r, 4).Value = 0
}

# Update F3 formula to sum the whole D2:D16 range instead of just D2
$ws.Range("F3").Formula = "=SUM(D2:D16)/F2"

# Update the active selection to match the author's final cursor position
$ws.Activate()
$ws.Range("I11").Select()
